$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 842.3333
$ws.Range("I2").Value = 842.3333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 842.3333
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -729.3333
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3284.7693
$ws.Range("J64").Value = 3311.111
$ws.Range("L64").Value = 3311.111
$ws.Range("N64").Value = -3807.111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3284.7693
$ws.Range("J67").Value = 3311.111
$ws.Range("L67").Value = 3311.111
$ws.Range("N67").Value = -5027.111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1066.6666
$ws.Range("I70").Value = 1200
$ws.Range("K70").Value = 3600
$ws.Range("M70").Value = -3330

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1066.6666
$ws.Range("I73").Value = 1200
$ws.Range("K73").Value = 3600
$ws.Range("M73").Value = -2664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I86").Value = 4986
$ws.Range("J86").Value = 4993.5
$ws.Range("K86").Value = 4986
$ws.Range("L86").Value = 4993.5
$ws.Range("M86").Value = -3863
$ws.Range("N86").Value = -7239.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I89").Value = 4986
$ws.Range("J89").Value = 4993.5
$ws.Range("K89").Value = 24930
$ws.Range("L89").Value = 24967.5
$ws.Range("M89").Value = -19314
$ws.Range("N89").Value = -36199.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 717.8823
$ws.Range("I98").Value = 820.6429000000001
$ws.Range("J98").Value = 238.33333
$ws.Range("K98").Value = 820.6429000000001
$ws.Range("L98").Value = 238.33333
$ws.Range("M98").Value = 677.3570999999999
$ws.Range("N98").Value = -3234.33333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 717.8823
$ws.Range("I122").Value = 820.6429000000001
$ws.Range("J122").Value = 238.33333
$ws.Range("K122").Value = 2461.9287
$ws.Range("L122").Value = 714.99999
$ws.Range("M122").Value = -11.92870000000039
$ws.Range("N122").Value = -5614.99999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2236.7144
$ws.Range("I2").Value = 2040
$ws.Range("J2").Value = 2499
$ws.Range("K2").Value = 2040
$ws.Range("L2").Value = 2499
$ws.Range("M2").Value = -1927
$ws.Range("N2").Value = -2725

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 49999
$ws.Range("J62").Value = 49999
$ws.Range("L62").Value = 49999
$ws.Range("N62").Value = -51247

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1000
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -314
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 49999
$ws.Range("J65").Value = 49999
$ws.Range("L65").Value = 149997
$ws.Range("N65").Value = -156237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1000
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 5000
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1392.6578
$ws.Range("I74").Value = 1003.9032
$ws.Range("K74").Value = 1003.9032
$ws.Range("M74").Value = -129.9032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1392.6578
$ws.Range("I77").Value = 1003.9032
$ws.Range("K77").Value = 5019.516
$ws.Range("M77").Value = -651.5159999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2236.7144
$ws.Range("I116").Value = 2040
$ws.Range("J116").Value = 2499
$ws.Range("K116").Value = 2040
$ws.Range("L116").Value = 2499
$ws.Range("M116").Value = 254
$ws.Range("N116").Value = -7087

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3064.3572
$ws.Range("I132").Value = 2363.25
$ws.Range("K132").Value = 7089.75
$ws.Range("M132").Value = -4559.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2236.7144
$ws.Range("I3").Value = 2040
$ws.Range("J3").Value = 2499
$ws.Range("K3").Value = 2040
$ws.Range("L3").Value = 2499
$ws.Range("M3").Value = -1926
$ws.Range("N3").Value = -2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2927.125
$ws.Range("I58").Value = 2898.1667
$ws.Range("K58").Value = 2898.1667
$ws.Range("M58").Value = -2695.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 5000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 5000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4632
$ws.Range("I132").Value = 4448.5
$ws.Range("K132").Value = 13345.5
$ws.Range("M132").Value = -10815.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2927.125
$ws.Range("I136").Value = 2898.1667
$ws.Range("K136").Value = 8694.500100000001
$ws.Range("M136").Value = -6144.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = -14888
$ws.Range("N3").Value = -15224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3464.889
$ws.Range("J113").Value = 3398.625
$ws.Range("L113").Value = 10195.875
$ws.Range("N113").Value = -14535.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 30
$ws.Range("I136").Value = 30
$ws.Range("K136").Value = 90
$ws.Range("M136").Value = 5010

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4749.5
$ws.Range("J80").Value = 7000
$ws.Range("L80").Value = 7000
$ws.Range("N80").Value = -8996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4749.5
$ws.Range("J83").Value = 7000
$ws.Range("L83").Value = 35000
$ws.Range("N83").Value = -44984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4048.652
$ws.Range("I132").Value = 3901.158
$ws.Range("K132").Value = 11703.474
$ws.Range("M132").Value = -9173.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5781.2
$ws.Range("I132").Value = 5636
$ws.Range("K132").Value = 16908
$ws.Range("M132").Value = -14378

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4334.36
$ws.Range("I136").Value = 4476.7896
$ws.Range("J136").Value = 3883.3333
$ws.Range("K136").Value = 13430.3688
$ws.Range("L136").Value = 11649.9999
$ws.Range("M136").Value = -10880.3688
$ws.Range("N136").Value = -16749.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 125896.25
$ws.Range("I62").Value = 167262
$ws.Range("J62").Value = 1799
$ws.Range("K62").Value = 167262
$ws.Range("L62").Value = 1799
$ws.Range("M62").Value = -166638
$ws.Range("N62").Value = -3047

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 125896.25
$ws.Range("I65").Value = 167262
$ws.Range("J65").Value = 1799
$ws.Range("K65").Value = 836310
$ws.Range("L65").Value = 8995
$ws.Range("M65").Value = -833190
$ws.Range("N65").Value = -15235

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3677.5625
$ws.Range("I136").Value = 3485.5
$ws.Range("K136").Value = 10456.5
$ws.Range("M136").Value = -7906.5
